$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "John Doe"
$ws.Range("B3").Value = "thisisnotapassword"

$ws.Range("B3").Select()
